$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D11").Value = -7.733999999999993
$ws.Range("A12").Value = -21.40049999999998
$ws.Range("D23").Value = -8.408999999999999
$ws.Range("A27").Value = -21.9648
$ws.Range("D28").Value = -8.200599999999994
$ws.Range("A32").Value = -21.53910000000001
$ws.Range("D32").Value = -7.037699999999993
$ws.Range("D34").Value = -8.130200000000004
$ws.Range("A36").Value = -19.7687
$ws.Range("A38").Value = -19.48439999999998
$ws.Range("D42").Value = -8.914899999999996
$ws.Range("A46").Value = -21.9179
$ws.Range("D49").Value = -8.162400000000007
$ws.Range("A54").Value = -21.8943
$ws.Range("D54").Value = -7.795699999999996
$ws.Range("A55").Value = -22.3004
$ws.Range("A56").Value = -22.13800000000001
$ws.Range("A67").Value = -21.45219999999997
$ws.Range("A69").Value = -21.50349999999998
$ws.Range("A72").Value = -21.7871
$ws.Range("D78").Value = -8.072199999999999
$ws.Range("D80").Value = -7.874999999999996
$ws.Range("A83").Value = -21.63489999999998
$ws.Range("A86").Value = -21.88890000000001
$ws.Range("A91").Value = -20.38419999999998
$ws.Range("A93").Value = -21.45870000000001
$ws.Range("D97").Value = -8.460299999999991
$ws.Range("A99").Value = -21.83610000000001
$ws.Range("D99").Value = -7.659499999999996
$ws.Range("D101").Value = -7.941599999999991
$ws.Range("A104").Value = -21.5086
